# Updated Argent prices in Excel
# Append a new row (2025-04-11) to each price-history sheet, carrying
# forward the latest known price for that series.

$wb = $excel.ActiveWorkbook

$newDate = "2025-04-11"

# Sheet name -> price value to write into column B of the new row.
$updates = [ordered]@{
    "N-Dense"                   = "40"
    "N-Type"                    = "41.5"
    "N-type Wafer"              = "1.28"
    "Cell Topcon 183mm"         = "0.303"
    "Module Topcon 183mm"       = "0.1"
    "Silver Rear_side"          = "5,187"
    "Silver Busbar front-side"  = "7,765"
    "Silver finger front-side"  = "7,815"
    "USD_CNY"                   = "7.3659"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the first empty row after the existing data in column A.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    # Prefix with an apostrophe so Excel stores the date/number as literal
    # text (matching the existing column formatting), then clear the
    # resulting "quote prefix" cell format so no stray style lingers.
    $ws.Cells.Item($newRow, 1).Value = "'" + $newDate
    $ws.Cells.Item($newRow, 2).Value = "'" + $updates[$sheetName]
    $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 2)).ClearFormats()
}
